# Common: Added an ability to remove a vape
# Adds six new translation rows (keys + Czech values) to the "Import"
# sheet of the translations fixture workbook, mirroring the existing
# cs/key/value layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting (style incl. wrapText) of the last existing data
# row down onto the new rows before filling in their content, so the
# new cells pick up the same cell style ("import", wrapText) as the
# rest of the table.
$ws.Range("A553:C553").Copy()
$ws.Range("A554:C559").PasteSpecial(-4122)

$ws.Range("A554").Value = "cs"
$ws.Range("B554").Value = "lab.vape.button.delete.confirm"
$ws.Range("C554").Value = "Opravdu si přejete odstranit vybraný vape? Touto akcí ovlivníte statistiky pro buildy, atomizéry a další, které jsou postavené na datech o vapování."

$ws.Range("A555").Value = "cs"
$ws.Range("B555").Value = "lab.vape.button.delete.confirm.ok"
$ws.Range("C555").Value = "Odstranit vape"

$ws.Range("A556").Value = "cs"
$ws.Range("B556").Value = "common.cancel"
$ws.Range("C556").Value = "Zrušit"

$ws.Range("A557").Value = "cs"
$ws.Range("B557").Value = "lab.vape.button.delete"
$ws.Range("C557").Value = "Odstranit vape"

$ws.Range("A558").Value = "cs"
$ws.Range("B558").Value = "lab.vape.button.delete.confirm.title"
$ws.Range("C558").Value = "Odstranit vape"

$ws.Range("A559").Value = "cs"
$ws.Range("B559").Value = "lab.vape.deleted.success"
$ws.Range("C559").Value = "Vape byl úspěšně odstraněn."

# Row 554 wraps onto two lines because of the long confirmation text,
# so it needs the taller row height (matches the other multi-line rows
# in the sheet).
$ws.Range("A554").EntireRow.RowHeight = 26.25

# Reflect the new last-edited cell in the sheet view, like Excel would
# after typing the new rows in.
$ws.Activate()
$ws.Range("B554").Select()
